$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new rows at the top of this data block (rows 342-344),
# shifting the existing rows 342-435 down to 345-438.
$ws.Rows("342:344").Insert()

# New row 342 (Acelga, Extra)
$ws.Range("A342").Value = 9
$ws.Range("B342").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C342").Value = "Metropolitana"
$ws.Range("D342").Value = 44551
$ws.Range("E342").Value = 13
$ws.Range("F342").Value = 100112009
$ws.Range("G342").Value = "Acelga"
$ws.Range("H342").Value = "Sin especificar"
$ws.Range("I342").Value = "Extra"
$ws.Range("J342").Value = 43
$ws.Range("K342").Value = 15000
$ws.Range("L342").Value = 16000
$ws.Range("M342").Value = 15512
$ws.Range("N342").Value = "$/docena de atados"
$ws.Range("O342").Value = "Región Metropolitana"
$ws.Range("P342").Value = 5171
$ws.Range("Q342").Value = 3
$ws.Range("R342").Value = "Hortaliza"

# New row 343 (Acelga, Primera)
$ws.Range("A343").Value = 9
$ws.Range("B343").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C343").Value = "Metropolitana"
$ws.Range("D343").Value = 44551
$ws.Range("E343").Value = 13
$ws.Range("F343").Value = 100112009
$ws.Range("G343").Value = "Acelga"
$ws.Range("H343").Value = "Sin especificar"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 61
$ws.Range("K343").Value = 13000
$ws.Range("L343").Value = 14000
$ws.Range("M343").Value = 13508
$ws.Range("N343").Value = "$/docena de atados"
$ws.Range("O343").Value = "Región Metropolitana"
$ws.Range("P343").Value = 4503
$ws.Range("Q343").Value = 3
$ws.Range("R343").Value = "Hortaliza"

# New row 344 (Acelga, Segunda)
$ws.Range("A344").Value = 9
$ws.Range("B344").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C344").Value = "Metropolitana"
$ws.Range("D344").Value = 44551
$ws.Range("E344").Value = 13
$ws.Range("F344").Value = 100112009
$ws.Range("G344").Value = "Acelga"
$ws.Range("H344").Value = "Sin especificar"
$ws.Range("I344").Value = "Segunda"
$ws.Range("J344").Value = 34
$ws.Range("K344").Value = 11000
$ws.Range("L344").Value = 12000
$ws.Range("M344").Value = 11500
$ws.Range("N344").Value = "$/docena de atados"
$ws.Range("O344").Value = "Región Metropolitana"
$ws.Range("P344").Value = 3833
$ws.Range("Q344").Value = 3
$ws.Range("R344").Value = "Hortaliza"
